# Fix the "property_category" column for the single-row "building" (建物)
# and "car" (汽車) property sheets: both were left over as "land" from a
# copy/paste of the original 土地 (Land) sheet. Set them to their correct
# category values.
#
# #5: property aircraft done

$wb = $excel.ActiveWorkbook

# 建物 (Building) sheet - property_category column is I, data row 2
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"

# 汽車 (Car) sheet - property_category column is H, data row 2
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
